$wb = $excel.ActiveWorkbook

# --- "Dest. Signs" sheet: update row 2, remove row 3 ---
$ws = $wb.Worksheets.Item("Dest. Signs")
$ws.Range("A2").ClearContents()
$ws.Range("B2").Value = "Emergency button does not work"
$ws.Range("C2").ClearContents()
$ws.Rows.Item(3).Delete()

# --- Other sheets: remove the single data row (row 2) ---
$sheetNames = @("Mirrors", "Other", "Zonar", "Stop Request", "Radio & PA")
foreach ($name in $sheetNames) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Rows.Item(2).Delete()
}
